$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
Write-Host $ws.Name
Write-Host $ws.Range("A6").Value
